# Error Calculations and Plots
# Applies the updated imputed/removed values to the missing_data sheet,
# and drops two rows (RM 232 and SC 92) that were removed from the
# combination_3_ABCDF / ACDF / 20 / seed5 sample.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Cell edits for rows 2-25 (row numbers unaffected by the later
#        row deletions, since both deleted rows are below row 25) ---
$ws.Range("D2").Value = -13.5    # RM 2
$ws.Range("F2").Value = 18.03    # RM 2
$ws.Range("E3").Value = ""       # RM 8
$ws.Range("E4").Value = -6.4     # RM 9
$ws.Range("D6").Value = ""       # RM 21
$ws.Range("E8").Value = ""       # RM 38
$ws.Range("F8").Value = ""       # RM 38
$ws.Range("E9").Value = ""       # RM 42
$ws.Range("F9").Value = ""       # RM 42
$ws.Range("D12").Value = -14.1   # RM 81
$ws.Range("F13").Value = 17.1    # RM 88
$ws.Range("D14").Value = ""      # RM 90
$ws.Range("E15").Value = -8.4    # RM 95
$ws.Range("E18").Value = -8.5    # RM 120
$ws.Range("E19").Value = ""      # RM 125
$ws.Range("F19").Value = ""      # RM 125
$ws.Range("D20").Value = -14     # RM 134
$ws.Range("D21").Value = -14.3   # RM 135
$ws.Range("E22").Value = ""      # RM 138
$ws.Range("D23").Value = ""      # RM 140
$ws.Range("E23").Value = -7      # RM 140
$ws.Range("D24").Value = ""      # RM 142a
$ws.Range("E25").Value = -7.1    # RM 145
$ws.Range("F25").Value = 16.6    # RM 145

# --- 2) Remove the "RM 232" and "SC 92" rows entirely (entire-row
#        delete, shifting everything below up) ---
$ws.Rows.Item(28).Delete()   # SC 92  (delete higher row index first)
$ws.Rows.Item(26).Delete()   # RM 232

# --- 3) Cell edits for the rows that now sit at their final positions
#        after the two deletions above ---
$ws.Range("B26").Value = -20.2   # SC 5
$ws.Range("B27").Value = ""      # SC 101
$ws.Range("E27").Value = ""      # SC 101
$ws.Range("B28").Value = ""      # SC 105
$ws.Range("F28").Value = 17.44   # SC 105
$ws.Range("B29").Value = -19.5   # SC 119
$ws.Range("B30").Value = -19.7   # SC 120
$ws.Range("B31").Value = ""      # SC 132
$ws.Range("D31").Value = -13.7   # SC 132
$ws.Range("F31").Value = ""      # SC 132
$ws.Range("B32").Value = ""      # SC 193
$ws.Range("F32").Value = 17.39   # SC 193
$ws.Range("D33").Value = -14.1   # SC 232
